$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.862.50'
$ws.Range('E2').Value = '  -1.10%  '
$ws.Range('D3').Value = '2.998.17'
$ws.Range('E3').Value = '  -0.48%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.84'
$ws.Range('E5').Value = '  +1.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.55'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '2.998.11'
$ws.Range('E8').Value = '  -0.35%  '
$ws.Range('E9').Value = '  -2.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.21'
$ws.Range('E10').Value = '  +6.94%  '
$ws.Range('E11').Value = '  -0.45%  '
$ws.Range('E12').Value = '  -1.02%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.20'
$ws.Range('E14').Value = '  -1.80%  '
$ws.Range('E15').Value = '  +2.69%  '
$ws.Range('D16').Value = '3.494.12'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.98'
$ws.Range('E17').Value = '  -1.81%  '
$ws.Range('D18').Value = '61.888.78'
$ws.Range('E18').Value = '  -1.00%  '
$ws.Range('D19').Value = '3.005.49'
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '447.28'
$ws.Range('E20').Value = '  -2.75%  '
$ws.Range('E21').Value = '  +0.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.684'
$ws.Range('E22').Value = '  -0.97%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.36'
$ws.Range('E23').Value = '  -1.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.01'
$ws.Range('E24').Value = '  +0.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.91'
$ws.Range('E25').Value = '  +8.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.21'
$ws.Range('E26').Value = '  -0.36%  '
$ws.Range('E27').Value = '  -2.23%  '
$ws.Range('E28').Value = '  +0.14%  '
$ws.Range('E29').Value = '  +2.66%  '
$ws.Range('B30').Value = 'FirstDigitalUSD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.15%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.21'
$ws.Range('E31').Value = '  +1.91%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.08'
$ws.Range('E32').Value = '  -0.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.40'
$ws.Range('E33').Value = '  -2.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.110'
$ws.Range('E34').Value = '  +0.75%  '
$ws.Range('D35').Value = '0.0₃0842'
$ws.Range('E35').Value = '  +3.99%  '
$ws.Range('E36').Value = '  -1.13%  '
$ws.Range('E37').Value = '  +0.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '50.10'
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('E39').Value = '  -4.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.95'
$ws.Range('E40').Value = '  -2.12%  '
$ws.Range('E41').Value = '  +4.04%  '
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '395.84'
$ws.Range('E43').Value = '  +0.86%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.51'
$ws.Range('E44').Value = '  +8.87%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.278'
$ws.Range('E45').Value = '  +3.32%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0349'
$ws.Range('E46').Value = '  -2.68%  '
$ws.Range('D47').Value = '2.709.55'
$ws.Range('E47').Value = '  -0.99%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '132.39'
$ws.Range('E48').Value = '  +2.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.16'
$ws.Range('E50').Value = '  -1.91%  '
$ws.Range('E51').Value = '  -1.87%  '
